# Saldo_guide.xlsx update: refresh the daily client-balance export.
# - Reference date (column G, "Dt. Referencia") rolls from 2024-11-18 to 2024-11-19
#   for every data row (2..274), matching the new export timestamp.
# - The worksheet/tab name is refreshed to the new export's timestamp.
# - A handful of accounts got revised Saldo Previsto / Vl. Total figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export run (2024-11-19 09:51:22).
$ws.Name = "IClientBalance-20241119-095122-"

# Bump the "Dt. Referencia" date serial for every data row (2-274) from 45614 to 45615.
$ws.Range("G2:G274").Value = 45615

# Revised balances for a handful of accounts (Saldo Previsto in column E,
# Vl. Total in column H move together).
$ws.Range("E6").Value = 999.99
$ws.Range("H6").Value = 999.99

$ws.Range("E33").Value = 614.85
$ws.Range("H33").Value = 614.85

$ws.Range("E49").Value = 419.84
$ws.Range("H49").Value = 419.84

$ws.Range("E107").Value = 12.13
$ws.Range("H107").Value = 12.13

$ws.Range("E255").Value = 219.59
$ws.Range("H255").Value = 219.59
